$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old label text -> new label text (x-axis bucket labels)
$map = @(
    @("Less than 200", "Less than 10K"),
    @("200 and less than  400", "10K and less than  21K"),
    @("400 and less than 600", "20K and less than 31K"),
    @("600 and less than 800", "31K and less than 42K"),
    @("800 and less than 1000", "42K and less than 52K"),
    @("1000 and less than 1200", "52K and less than 62K"),
    @("1200 and less than 1400", "62K and less than 73K"),
    @("1400 and less than 1600", "73K and less than 84K"),
    @("1600 and less than 1800", "84K and less than 94K"),
    @("1800 and less than 2000", "94K and less than 100K"),
    @("Above 2000", "Above 100K")
)

$xlWhole = 1
$xlByRows = 1

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $ws.Cells.Replace($old, $new, $xlWhole, $xlByRows, $false, $false, $false)
}
